# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts and one "最低票价" (G column)
# fix across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 251
$ws.Range("F3").Value  = 589
$ws.Range("F6").Value  = 2784
$ws.Range("F8").Value  = 49
$ws.Range("F10").Value = 362
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 301
$ws.Range("F14").Value = 5820
$ws.Range("F16").Value = 1030
$ws.Range("F17").Value = 75
$ws.Range("F20").Value = 489
$ws.Range("F22").Value = 1267
$ws.Range("F23").Value = 77
$ws.Range("F24").Value = 16
$ws.Range("F25").Value = 2029
$ws.Range("F26").Value = 143
$ws.Range("F27").Value = 340
$ws.Range("F29").Value = 3210

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 17
$ws.Range("F25").Value = 4039
$ws.Range("F32").Value = 190
$ws.Range("G33").Value = 480

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value  = 2580
$ws.Range("F9").Value  = 1456
$ws.Range("F10").Value = 408
$ws.Range("F13").Value = 592

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 2580
$ws.Range("F7").Value  = 1456
$ws.Range("F8").Value  = 408
$ws.Range("F11").Value = 251
$ws.Range("F12").Value = 589
$ws.Range("F14").Value = 2784
$ws.Range("F15").Value = 49
$ws.Range("F17").Value = 592
$ws.Range("F18").Value = 362
$ws.Range("F21").Value = 301
$ws.Range("F23").Value = 5820
$ws.Range("F25").Value = 1030
$ws.Range("F27").Value = 75
$ws.Range("F30").Value = 489
$ws.Range("F38").Value = 1267
$ws.Range("F41").Value = 16
$ws.Range("F43").Value = 2029
$ws.Range("F44").Value = 190
$ws.Range("G45").Value = 480
$ws.Range("F46").Value = 143
$ws.Range("F47").Value = 340
$ws.Range("F49").Value = 3210

$wb.Save()
